# Weekly update: insert a new price record row for
# "Hortaliza, Femacal de La Calera - Apio" at the top of the
# date-ordered block (row 175), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 175; this shifts rows 175-239 down to 176-240
# and preserves formatting (e.g. the date style on column D).
$ws.Rows("175:175").Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(175, 1).Value = 3
$ws.Cells.Item(175, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(175, 3).Value = "Coquimbo"
$ws.Cells.Item(175, 4).Value = 44468
$ws.Cells.Item(175, 5).Value = 5
$ws.Cells.Item(175, 6).Value = 100112017
$ws.Cells.Item(175, 7).Value = "Apio"
$ws.Cells.Item(175, 8).Value = "Americana (o)"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 240
$ws.Cells.Item(175, 11).Value = 8500
$ws.Cells.Item(175, 12).Value = 9000
$ws.Cells.Item(175, 13).Value = 8812
$ws.Cells.Item(175, 14).Value = "`$/docena de matas"
$ws.Cells.Item(175, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(175, 16).Value = 1469
$ws.Cells.Item(175, 17).Value = 6
$ws.Cells.Item(175, 18).Value = "Hortaliza"
